$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2830809886899317
$ws.Range("C2").Value = 0.05228370174917529
$ws.Range("D2").Value = 0.03306337985746666
$ws.Range("E2").Value = 0.1652192017219107
$ws.Range("F2").Value = 0.8129378803801686
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("K2").Value = 0.2587585588637182
$ws.Range("M2").Value = 0.2176760785492249
$ws.Range("O2").Value = 2.832204797704549
$ws.Range("B3").Value = 0.2507523373225808
$ws.Range("C3").Value = 0.04908269789008557
$ws.Range("D3").Value = 0.03135415875789249
$ws.Range("E3").Value = 0.1541697129981614
$ws.Range("F3").Value = 0.8105211648305612
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("K3").Value = 0.2259095237996718
$ws.Range("M3").Value = 0.1957731359861583
$ws.Range("O3").Value = 2.838189110714126
$ws.Range("B4").Value = 0.2309166325481726
$ws.Range("C4").Value = 0.04710047589719579
$ws.Range("D4").Value = 0.03029392869230918
$ws.Range("E4").Value = 0.1474994569335877
$ws.Range("F4").Value = 0.8095404201753809
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("K4").Value = 0.2056964118789892
$ws.Range("M4").Value = 0.1823977817766789
$ws.Range("O4").Value = 2.843528099202757
$ws.Range("B5").Value = 0.2228373722532808
$ws.Range("C5").Value = 0.04628852206164424
$ws.Range("D5").Value = 0.02985919589664832
$ws.Range("E5").Value = 0.1448097466896172
$ws.Range("F5").Value = 0.8092671998140588
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("K5").Value = 0.1974487944272028
$ws.Range("M5").Value = 0.1769655473930953
$ws.Range("O5").Value = 2.846122018784769
$ws.Range("B6").Value = 0.2214960669334687
$ws.Range("C6").Value = 0.04615344634317609
$ws.Range("D6").Value = 0.02978684771176177
$ws.Range("E6").Value = 0.1443648351531053
$ws.Range("F6").Value = 0.8092294663289934
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("K6").Value = 0.1960786539934958
$ws.Range("M6").Value = 0.1760646344298635
$ws.Range("O6").Value = 2.846577987859433
$ws.Range("B7").Value = 0.2308076563241457
$ws.Range("C7").Value = 0.04708954248432917
$ws.Range("D7").Value = 0.03028807655273624
$ws.Range("E7").Value = 0.1474630675859245
$ws.Range("F7").Value = 0.8095362235761954
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("K7").Value = 0.2055852240915499
$ws.Range("M7").Value = 0.1823244466210028
$ws.Range("O7").Value = 2.843561388733889
$ws.Range("B8").Value = 0.2719313333299169
$ws.Range("C8").Value = 0.05118350414022643
$ws.Range("D8").Value = 0.03247628878722253
$ws.Range("E8").Value = 0.1613854949031506
$ws.Range("F8").Value = 0.8120001246013402
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("K8").Value = 0.2474414878395237
$ws.Range("M8").Value = 0.210108719840747
$ws.Range("O8").Value = 2.833922492611549
$ws.Range("B9").Value = 0.3526751361911806
$ws.Range("C9").Value = 0.05907710337511674
$ws.Range("D9").Value = 0.03668104365515035
$ws.Range("E9").Value = 0.1896054350021572
$ws.Range("F9").Value = 0.8208287424903773
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("K9").Value = 0.3291626809500769
$ws.Range("M9").Value = 0.2651795853471128
$ws.Range("O9").Value = 2.828246759254682
$ws.Range("B10").Value = 0.4120481324848413
$ws.Range("C10").Value = 0.0647931022357966
$ws.Range("D10").Value = 0.03971667308567817
$ws.Range("E10").Value = 0.2109173303032392
$ws.Range("F10").Value = 0.829760513062908
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("K10").Value = 0.3889738941034295
$ws.Range("M10").Value = 0.3060093878631349
$ws.Range("O10").Value = 2.832169520643475
$ws.Range("B11").Value = 0.4390675447412207
$ws.Range("C11").Value = 0.06737509638203676
$ws.Range("D11").Value = 0.0410858313152076
$ws.Range("E11").Value = 0.2207425410970387
$ws.Range("F11").Value = 0.8343568068720941
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("K11").Value = 0.4161319071566254
$ws.Range("M11").Value = 0.3246669330915068
$ws.Range("O11").Value = 2.835718031894146
$ws.Range("B12").Value = 0.4493002885963904
$ws.Range("C12").Value = 0.06835017487489381
$ws.Range("D12").Value = 0.04160258225415703
$ws.Range("E12").Value = 0.2244821261929957
$ws.Range("F12").Value = 0.8361740954580625
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("K12").Value = 0.4264084038764793
$ws.Range("M12").Value = 0.3317442548471945
$ws.Range("O12").Value = 2.83731588945389
$ws.Range("B13").Value = 0.447096442612235
$ws.Range("C13").Value = 0.0681402935614841
$ws.Range("D13").Value = 0.04149136761240158
$ws.Range("E13").Value = 0.2236758910059891
$ws.Range("F13").Value = 0.835779293809253
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("K13").Value = 0.4241955232469081
$ws.Range("M13").Value = 0.330219487313272
$ws.Range("O13").Value = 2.836960453352702
$ws.Range("B14").Value = 0.4399093794446856
$ws.Range("C14").Value = 0.06745537036813687
$ws.Range("D14").Value = 0.04112837934686553
$ws.Range("E14").Value = 0.2210498169676853
$ws.Range("F14").Value = 0.8345047771042999
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("K14").Value = 0.4169775164228611
$ws.Range("M14").Value = 0.3252489455074041
$ws.Range("O14").Value = 2.835844393224164
$ws.Range("B15").Value = 0.4355072241566518
$ws.Range("C15").Value = 0.06703548657692693
$ws.Range("D15").Value = 0.04090581379371372
$ws.Range("E15").Value = 0.2194437514313776
$ws.Range("F15").Value = 0.8337340997446319
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("K15").Value = 0.4125552701514437
$ws.Range("M15").Value = 0.3222059244043507
$ws.Range("O15").Value = 2.835193880512804
$ws.Range("B16").Value = 0.4102825234811291
$ws.Range("C16").Value = 0.06462399204060887
$ws.Range("D16").Value = 0.03962695625553891
$ws.Range("E16").Value = 0.2102778735239355
$ws.Range("F16").Value = 0.8294708721264783
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("K16").Value = 0.387198005086816
$ws.Range("M16").Value = 0.3047917656287282
$ws.Range("O16").Value = 2.831973152188937
$ws.Range("B17").Value = 0.3948103437952
$ws.Range("C17").Value = 0.06313991490046078
$ws.Range("D17").Value = 0.03883938548380428
$ws.Range("E17").Value = 0.2046884613820552
$ws.Range("F17").Value = 0.8269921536285665
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("K17").Value = 0.3716289448530858
$ws.Range("M17").Value = 0.2941302615771662
$ws.Range("O17").Value = 2.830449468906806
$ws.Range("B18").Value = 0.3859121526270144
$ws.Range("C18").Value = 0.06228459848867374
$ws.Range("D18").Value = 0.03838528989558654
$ws.Range("E18").Value = 0.2014858309370595
$ws.Range("F18").Value = 0.8256166404156033
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("K18").Value = 0.3626693126244902
$ws.Range("M18").Value = 0.2880059367448027
$ws.Range("O18").Value = 2.829739098452364
$ws.Range("B19").Value = 0.3828995617446935
$ws.Range("C19").Value = 0.06199470991354872
$ws.Range("D19").Value = 0.03823135187441551
$ws.Range("E19").Value = 0.2004035718830863
$ws.Range("F19").Value = 0.8251595306883388
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("K19").Value = 0.3596349367105347
$ws.Range("M19").Value = 0.2859336990883747
$ws.Range("O19").Value = 2.829527078247423
$ws.Range("B20").Value = 0.3964572834730973
$ws.Range("C20").Value = 0.06329807520505426
$ws.Range("D20").Value = 0.0389233383760228
$ws.Range("E20").Value = 0.2052821940398672
$ws.Range("F20").Value = 0.827250823260016
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("K20").Value = 0.3732867893955074
$ws.Range("M20").Value = 0.2952643798721937
$ws.Range("O20").Value = 2.830594482727719
$ws.Range("B21").Value = 0.4420203680754184
$ws.Range("C21").Value = 0.06765662162811736
$ws.Range("D21").Value = 0.04123504468868333
$ws.Range("E21").Value = 0.2218206411590202
$ws.Range("F21").Value = 0.8348770491944322
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("K21").Value = 0.4190978309206628
$ws.Range("M21").Value = 0.3267085849058446
$ws.Range("O21").Value = 2.836165307007349
$ws.Range("B22").Value = 0.4718045780902571
$ws.Range("C22").Value = 0.07048962052712682
$ws.Range("D22").Value = 0.04273584266793762
$ws.Range("E22").Value = 0.2327402831699743
$ws.Range("F22").Value = 0.8403087597599779
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("K22").Value = 0.4489931986003342
$ws.Range("M22").Value = 0.3473298277707642
$ws.Range("O22").Value = 2.84128753696649
$ws.Range("B23").Value = 0.4559077607529787
$ws.Range("C23").Value = 0.06897903396023253
$ws.Range("D23").Value = 0.04193576594881421
$ws.Range("E23").Value = 0.2269020425681276
$ws.Range("F23").Value = 0.837368771610457
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("K23").Value = 0.4330417102545709
$ws.Range("M23").Value = 0.3363174008265517
$ws.Range("O23").Value = 2.838418012816817
$ws.Range("B24").Value = 0.3957127114505852
$ws.Range("C24").Value = 0.06322657747271876
$ws.Range("D24").Value = 0.03888538735424163
$ws.Range("E24").Value = 0.2050137338955835
$ws.Range("F24").Value = 0.8271337244265879
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("K24").Value = 0.3725373051731822
$ws.Range("M24").Value = 0.2947516290638603
$ws.Range("O24").Value = 2.830528406120379
$ws.Range("B25").Value = 0.3308221618610219
$ws.Range("C25").Value = 0.05695623117003379
$ws.Range("D25").Value = 0.03555288947520552
$ws.Range("E25").Value = 0.1818709411557364
$ws.Range("F25").Value = 0.8180115616049903
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("K25").Value = 0.3070944521543311
$ws.Range("M25").Value = 0.2502174838313849
$ws.Range("O25").Value = 2.828363172449144
